$d = $word.ActiveDocument

# Replacement 1
$found1 = $d.Content.Find.Execute("Ativação: 01/01/1996", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2021", 2)
if (-not $found1) { Write-Host "WARNING: replacement 1 not found" }

# Replacement 2
$found2 = $d.Content.Find.Execute("-Apresentar a evolução das condições geológicas da Terra que culminaram com os recursos naturais existentes hoje, com ênfase nas reservas de combustíveis fósseis, hídricos e de minérios e a conseqüente reserva de energia advinda dessas fontes naturais. A América do Sul e do Brasil, mais especificamente, são destacados.- Acompanhar, através da literatura atual, a escassez dos recursos em vista do mau uso e das alternativas para o aproveitamento racional dos recursos existentes, com alternativas na área de geração de energia.", $true, $false, $false, $false, $false, $true, 1, $false, "- Apresentar a evolução das condições geológicas da Terra que culminaram com os recursos naturais existentes hoje, com ênfase nas reservas de combustíveis fósseis, hídricos e de minérios e a conseqüente reserva de energia advinda dessas fontes naturais. A América do Sul e do Brasil, mais especificamente, são destacados.", 2)
if (-not $found2) { Write-Host "WARNING: replacement 2 not found" }

# Replacement 3
$found3 = $d.Content.Find.Execute("- Universo e Terra,- Desenvolvimento da Terra.- Combustíveis fósseis.- Recursos renováveis e biomassa.- Recursos minerais.- Matérias-primas da grande indústria química:metais ferrosos e não-ferrosos.- Recursos hídricos.", $true, $false, $false, $false, $false, $true, 1, $false, "- Desenvolvimento da Terra.- Recursos minerais.- Matérias-primas da grande indústria metalúrgica: metais ferrosos e não-ferrosos", 2)
if (-not $found3) { Write-Host "WARNING: replacement 3 not found" }

# Replacement 4
$found4 = $d.Content.Find.Execute("- Formação do Universo. - Formação do Sistema Solar.- Desenvolvimento da Terra.- Principais Eras Geológicas.- Petróleo.- Carvão e Gás Natural - Geração de  energia termelétrica.- Recursos Renováveis.- Biomassa ? Fontes alternativas de energia.- Matérias-primas para a grande indústria química.- Metais ferrosos.- Metais não-ferrosos.- Recursos hídricos ? Bacias hídricas.- Poluição das águas.- Escassez e reaproveitamento das águas.- Geração de energia elétrica.", $true, $false, $false, $false, $false, $true, 1, $false, "- Desenvolvimento da Terra. - Principais Eras Geológicas. - Matérias-primas para a grande indústria metalúrgica: metais ferrosos e metais não-ferrosos.", 2)
if (-not $found4) { Write-Host "WARNING: replacement 4 not found" }

# Replacement 5
$found5 = $d.Content.Find.Execute("- Schäfer, A., Fundamentos de Ecologia e Biogeografia de Águas Continentais, Ed. Universidade, Porto Alegre.- Abreu, S.F. Recursos Minerais do Brasil, Ed. Edgard Bluecher, 1973.- Carioca, J. O. B. and Arora, H.L. Biomassa-Fundamentos e Aplicações Tecnológicas, Universidade Federal do Ceará, 1984.- Fernandes, F.R.C. Quem é quem no Subsolo Brasileiro, MCT/CNPq, 1987- Petri, S. e Fúlfaro, V. Geologia do Brasil, EDUSP, 1983.- Revistas especializadas e fontes de informação multimídia as mais diversas, dado ao caráter dinâmico das informações sobre reservas minerais, geração de energia e recursos naturais em geral.", $true, $false, $false, $false, $false, $true, 1, $false, "- MILLER, Jr. G. T. “Ciência Ambiental”,  Editora: Thomson, 2006.- ABREU, S. F. “Recursos Minerais do Brasil”, Editora: Edgard Blücher, 1973.-  SKINNER, B. J. “Recursos Minerais da Terra”, Editora: Edgard Blücher, 1996.- WICANDER, R.; MONROE, J. S. “Fundamentos de Geologia”, Editora: Cengage Learning, 2009. - PRESS, F.; Siever, R.; Jordan, T.; Grotzinger, J. “Para Entender a Terra”, Editora: Bookman,  2006.- SCHÄFER, A. “Fundamentos de Ecologia e Biogeografia de Águas Continentais”, Editora: Universidade, Porto Alegre. - Revistas especializadas e sites, dado ao caráter dinâmico das informações sobre reservas minerais e recursos naturais em geral.", 2)
if (-not $found5) { Write-Host "WARNING: replacement 5 not found" }
